# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-35 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 2
    4  = 1
    5  = 9
    6  = 6
    7  = 2
    8  = 5
    9  = 6
    10 = 6
    11 = 7
    12 = 0
    13 = 6
    14 = 3
    15 = 8
    16 = 8
    17 = 4
    18 = 6
    19 = 7
    20 = 8
    21 = 6
    22 = 3
    23 = 3
    24 = 7
    25 = 4
    26 = 9
    27 = 6
    28 = 3
    29 = 2
    30 = 5
    31 = 4
    32 = 9
    33 = 2
    34 = 3
    35 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
